$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting existing rows 116:190 down to 117:191
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new data record
$ws.Cells.Item(116, 1).Value = 1
$ws.Cells.Item(116, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(116, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(116, 4).Value = 45029
$ws.Cells.Item(116, 5).Value = 15
$ws.Cells.Item(116, 6).Value = 100114001
$ws.Cells.Item(116, 7).Value = "Papa"
$ws.Cells.Item(116, 8).Value = "Red Lady"
$ws.Cells.Item(116, 9).Value = "1a (cosecha)"
$ws.Cells.Item(116, 10).Value = 1000
$ws.Cells.Item(116, 11).Value = 12000
$ws.Cells.Item(116, 12).Value = 13000
$ws.Cells.Item(116, 13).Value = 12500
$ws.Cells.Item(116, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(116, 15).Value = "Región del Bíobío"
$ws.Cells.Item(116, 16).Value = 500
$ws.Cells.Item(116, 17).Value = 25
$ws.Cells.Item(116, 18).Value = "Hortaliza"
